# "Add files via upload" -- refreshed copy of Henriette Wilford's keyword
# expertise review: update her self-reported level of expertise (L/M/H)
# for each keyword in column C of the "4 Henriette Wilford" sheet.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("4 Henriette Wilford")
$ws3.Range("C2").Value = "H"
$ws3.Range("C3").Value = "M"
$ws3.Range("C4").Value = "L"
$ws3.Range("C5").Value = "L"
$ws3.Range("C6").Value = "L"
$ws3.Range("C8").Value = ""
$ws3.Range("C10").Value = "M"
$ws3.Range("C11").Value = ""
$ws3.Range("C12").Value = "L"

# Leave the worksheets' prior (frozen-pane) selections as they were: touch
# row 10 (as the reviewer did while scanning the sheet) and then restore
# the original working selection in the scrollable pane.
$ws3.Activate()
$ws3.Range("A10:K10").Select()
$ws3.Range("C2:C101").Select()

$ws2 = $wb.Worksheets.Item("Conflicts of Interest")
$ws2.Activate()
$ws2.Range("A10:K10").Select()
$ws2.Range("F2:F26").Select()

# Restore the originally active sheet/tab.
$ws1 = $wb.Worksheets.Item("Expertise by Keywords - Instr.")
$ws1.Activate()
